$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 67 mirrors the layout of row 66 (the previous last data row).
# Copy A66's formatting (the bold/bordered date style) down to A67 first,
# then fill in the values for the new row.
$ws.Range("A66").Copy()
$ws.Range("A67").PasteSpecial(-4122)

$ws.Range("A67").Value = 45497
$ws.Range("B67").Value = 697.7737382400001
$ws.Range("C67").Value = 231.2928896935
$ws.Range("D67").Value = 0
$ws.Range("E67").Value = 0
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 121.66922805
$ws.Range("I67").Value = 261.1745245114
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 20.998543188008
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 131.492286616
$ws.Range("O67").Value = 58.395680213
$ws.Range("P67").Value = 0
$ws.Range("Q67").Value = 0.0000029328
$ws.Range("R67").Value = 0
$ws.Range("S67").Value = 0
$ws.Range("T67").Value = 0
$ws.Range("U67").Value = 346.4317811017268
$ws.Range("W67").Value = 0
$ws.Range("X67").Value = 0
$ws.Range("Y67").Value = 0
$ws.Range("Z67").Value = 253.878853961308
